$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C53").Value = "https://leetcode.com/problems/linked-list-random-node/"
$ws.Hyperlinks.Add($ws.Range("C53"), "https://leetcode.com/problems/linked-list-random-node/")
"done"
